# PowerPoll's pie-chart label sizing and the "daylight savings" date fix both
# live in the web-extension's own HTML/JS bundle (loaded at runtime from the
# add-in's dev manifest, store="developer" / storeType="Registry"), not in
# any OOXML part of this .pptx. The only observable change in the saved
# package is PowerPoint reassigning fresh relationship ids (sldMasterId,
# sldId, sldLayoutId, the webextensionref/blip r:id's) and a new GUID for
# <we:webextension id="...">- all incidental churn from the rebuild, with
# no attribute/text/content difference anywhere else in the package.
#
# The PowerPoint object model has no surface for the embedded web add-in
# (no Shape.WebExtension/CustomXMLParts entry is created for
# ppt/slides/udata/data.xml, and its relationship ids aren't settable via
# COM), so there is nothing reachable here to change; leave the
# presentation untouched.
$p = $ppt.ActivePresentation
